$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.627.19'
$ws.Range("E2").Value = '  -2.71%  '
$ws.Range("D3").Value = '1.666.60'
$ws.Range("E3").Value = '  -4.05%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.59'
$ws.Range("E5").Value = '  -2.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.510'
$ws.Range("E6").Value = '  -2.59%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.19'
$ws.Range("E8").Value = '  -0.59%  '
$ws.Range("E9").Value = '  -1.17%  '
$ws.Range("E10").Value = '  -2.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0879'
$ws.Range("E11").Value = '  -1.92%  '
$ws.Range("D12").Value = '1.902.39'
$ws.Range("E12").Value = '  -3.94%  '
$ws.Range("D13").Value = '1.660.57'
$ws.Range("E13").Value = '  -4.51%  '
$ws.Range("E14").Value = '  -3.44%  '
$ws.Range("E15").Value = '  +0.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.52'
$ws.Range("E16").Value = '  -1.94%  '
$ws.Range("D17").Value = '27.607.95'
$ws.Range("E17").Value = '  -2.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '242.26'
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("E19").Value = '  -3.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.71'
$ws.Range("E20").Value = '  -4.12%  '
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("E22").Value = '  -3.37%  '
$ws.Range("E23").Value = '  -3.55%  '
$ws.Range("E24").Value = '  -3.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.56'
$ws.Range("E25").Value = '  -1.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.25'
$ws.Range("E26").Value = '  -4.18%  '
$ws.Range("E27").Value = '  -1.99%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("E29").Value = '  -2.57%  '
$ws.Range("E30").Value = '  +1.90%  '
$ws.Range("E31").Value = '  -1.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.36'
$ws.Range("E32").Value = '  -2.70%  '
$ws.Range("D33").Value = '1.461.98'
$ws.Range("E33").Value = '  -2.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.14'
$ws.Range("E34").Value = '  -4.36%  '
$ws.Range("E35").Value = '  -4.95%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.931'
$ws.Range("E36").Value = '  -4.09%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.37'
$ws.Range("E37").Value = '  -1.30%  '
$ws.Range("E38").Value = '  -4.85%  '
$ws.Range("E39").Value = '  -2.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '70.05'
$ws.Range("E40").Value = '  -1.23%  '
$ws.Range("E41").Value = '  -4.24%  '
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.23'
$ws.Range("E43").Value = '  -3.53%  '
$ws.Range("E44").Value = '  -5.71%  '
$ws.Range("E45").Value = '  -1.39%  '
$ws.Range("D46").Value = '1.809.88'
$ws.Range("E46").Value = '  -3.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.75'
$ws.Range("E47").Value = '  +1.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.98'
$ws.Range("E48").Value = '  -2.42%  '
$ws.Range("D49").Value = '0.0₆0107'
$ws.Range("E49").Value = '  -4.60%  '
$ws.Range("E50").Value = '  -2.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.94'
$ws.Range("E51").Value = '  -4.10%  '
